$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trials")

$ws.Range("B402").Value = 14627
$ws.Range("C402").Value = 25.14061212539673
$ws.Range("B403").Value = 1233
$ws.Range("C403").Value = 0.6089184284210205
$ws.Range("B404").Value = 1972
$ws.Range("C404").Value = 1.199352502822876
$ws.Range("B405").Value = 2804
$ws.Range("C405").Value = 1.376344919204712
$ws.Range("B406").Value = 6613
$ws.Range("C406").Value = 6.687793731689453
$ws.Range("B407").Value = 6684
$ws.Range("C407").Value = 4.609642744064331
$ws.Range("B408").Value = 9944
$ws.Range("C408").Value = 11.15787434577942
$ws.Range("B409").Value = 4121
$ws.Range("C409").Value = 3.051645278930664
$ws.Range("B410").Value = 24470
$ws.Range("C410").Value = 60.00059390068054
$ws.Range("D410").Value = 0
$ws.Range("B411").Value = 1057
$ws.Range("C411").Value = 0.4389910697937012
$ws.Range("B412").Value = 4366
$ws.Range("C412").Value = 2.816287040710449
$ws.Range("B413").Value = 2701
$ws.Range("C413").Value = 1.417914390563965
$ws.Range("B414").Value = 11918
$ws.Range("C414").Value = 14.80842900276184
$ws.Range("B415").Value = 1961
$ws.Range("C415").Value = 0.9647173881530762
$ws.Range("D415").Value = 1
$ws.Range("B416").Value = 5814
$ws.Range("C416").Value = 4.525143384933472
$ws.Range("B417").Value = 14488
$ws.Range("C417").Value = 20.71092939376831
$ws.Range("B418").Value = 2825
$ws.Range("C418").Value = 1.45424485206604
$ws.Range("B419").Value = 3850
$ws.Range("C419").Value = 2.292064428329468
$ws.Range("B420").Value = 8822
$ws.Range("C420").Value = 8.336741924285889
$ws.Range("B421").Value = 3526
$ws.Range("C421").Value = 1.919137239456177
$ws.Range("B422").Value = 2662
$ws.Range("C422").Value = 1.337094068527222
$ws.Range("B423").Value = 3206
$ws.Range("C423").Value = 1.678504467010498
$ws.Range("B424").Value = 5380
$ws.Range("C424").Value = 3.761575222015381
$ws.Range("B425").Value = 2518
$ws.Range("C425").Value = 1.197699069976807
$ws.Range("B426").Value = 3245
$ws.Range("C426").Value = 1.545844793319702
$ws.Range("B427").Value = 1802
$ws.Range("C427").Value = 0.5446517467498779
$ws.Range("B428").Value = 7636
$ws.Range("C428").Value = 6.609802007675171
$ws.Range("B429").Value = 7211
$ws.Range("C429").Value = 5.916283130645752
$ws.Range("B430").Value = 1676
$ws.Range("C430").Value = 0.6909406185150146
$ws.Range("B431").Value = 6318
$ws.Range("C431").Value = 4.753081560134888
$ws.Range("B432").Value = 2829
$ws.Range("C432").Value = 1.048093795776367
$ws.Range("B433").Value = 6529
$ws.Range("C433").Value = 4.552799224853516
$ws.Range("B434").Value = 6207
$ws.Range("C434").Value = 4.709862947463989
$ws.Range("B435").Value = 8991
$ws.Range("C435").Value = 8.628479242324829
$ws.Range("B436").Value = 981
$ws.Range("C436").Value = 0.2667269706726074
$ws.Range("B437").Value = 4736
$ws.Range("C437").Value = 3.122681379318237
$ws.Range("B438").Value = 5044
$ws.Range("C438").Value = 3.305302858352661
$ws.Range("B439").Value = 2514
$ws.Range("C439").Value = 1.257573366165161
$ws.Range("B440").Value = 2494
$ws.Range("C440").Value = 1.149695158004761
$ws.Range("B441").Value = 17782
$ws.Range("C441").Value = 28.37506628036499
$ws.Range("D441").Value = 1
$ws.Range("B442").Value = 434
$ws.Range("C442").Value = 0.143179178237915
$ws.Range("B443").Value = 7247
$ws.Range("C443").Value = 7.52823805809021
$ws.Range("B444").Value = 1229
$ws.Range("C444").Value = 0.6102972030639648
$ws.Range("B445").Value = 10791
$ws.Range("C445").Value = 10.39438796043396
$ws.Range("B446").Value = 5086
$ws.Range("C446").Value = 3.460047960281372
$ws.Range("B447").Value = 4479
$ws.Range("C447").Value = 2.756176948547363
$ws.Range("B448").Value = 13583
$ws.Range("C448").Value = 16.96275520324707
$ws.Range("B449").Value = 12504
$ws.Range("C449").Value = 15.73604607582092
$ws.Range("B450").Value = 816
$ws.Range("C450").Value = 0.3293182849884033
$ws.Range("B451").Value = 1668
$ws.Range("C451").Value = 0.8283329010009766
$ws.Range("B452").Value = 2729
$ws.Range("C452").Value = 1.395344018936157
$ws.Range("D452").Value = 1
$ws.Range("B453").Value = 10816
$ws.Range("C453").Value = 12.46210789680481
$ws.Range("B454").Value = 2049
$ws.Range("C454").Value = 1.211089134216309
$ws.Range("B455").Value = 6526
$ws.Range("C455").Value = 5.146907567977905
$ws.Range("B456").Value = 10439
$ws.Range("C456").Value = 14.74250841140747
$ws.Range("B457").Value = 12344
$ws.Range("C457").Value = 15.48947548866272
$ws.Range("B458").Value = 9768
$ws.Range("C458").Value = 11.96901822090149
$ws.Range("B459").Value = 2324
$ws.Range("C459").Value = 1.431987285614014
$ws.Range("B460").Value = 6551
$ws.Range("C460").Value = 4.847487688064575
$ws.Range("B461").Value = 2811
$ws.Range("C461").Value = 1.417273998260498
$ws.Range("B462").Value = 5827
$ws.Range("C462").Value = 4.531213045120239
$ws.Range("B463").Value = 6892
$ws.Range("C463").Value = 6.076826572418213
$ws.Range("B464").Value = 4817
$ws.Range("C464").Value = 3.51691722869873
$ws.Range("B465").Value = 3302
$ws.Range("C465").Value = 2.01761531829834
$ws.Range("B466").Value = 7495
$ws.Range("C466").Value = 5.83779501914978
$ws.Range("B467").Value = 11124
$ws.Range("C467").Value = 14.93817353248596
$ws.Range("B468").Value = 8861
$ws.Range("C468").Value = 8.138381719589233
$ws.Range("D468").Value = 1
$ws.Range("B469").Value = 5687
$ws.Range("C469").Value = 3.718128204345703
$ws.Range("B470").Value = 7253
$ws.Range("C470").Value = 6.361289501190186
$ws.Range("B471").Value = 10042
$ws.Range("C471").Value = 10.73122453689575
$ws.Range("B472").Value = 969
$ws.Range("C472").Value = 0.4135811328887939
$ws.Range("B473").Value = 1907
$ws.Range("C473").Value = 0.8192036151885986
$ws.Range("B474").Value = 7147
$ws.Range("C474").Value = 6.35259485244751
$ws.Range("B475").Value = 8361
$ws.Range("C475").Value = 8.280720233917236
$ws.Range("B476").Value = 14313
$ws.Range("C476").Value = 19.11187982559204
$ws.Range("B477").Value = 14415
$ws.Range("C477").Value = 19.62747120857239
$ws.Range("B478").Value = 17818
$ws.Range("C478").Value = 33.17127847671509
$ws.Range("B479").Value = 6020
$ws.Range("C479").Value = 4.386992216110229
$ws.Range("B480").Value = 4193
$ws.Range("C480").Value = 2.675059795379639
$ws.Range("B481").Value = 8301
$ws.Range("C481").Value = 7.533198118209839
$ws.Range("B482").Value = 4576
$ws.Range("C482").Value = 2.953443288803101
$ws.Range("B483").Value = 996
$ws.Range("C483").Value = 0.3685545921325684
$ws.Range("B484").Value = 11864
$ws.Range("C484").Value = 14.77514719963074
$ws.Range("D484").Value = 1
$ws.Range("B485").Value = 5530
$ws.Range("C485").Value = 4.033508062362671
$ws.Range("B486").Value = 10476
$ws.Range("C486").Value = 12.53589534759521
$ws.Range("B487").Value = 8423
$ws.Range("C487").Value = 7.074377536773682
$ws.Range("B488").Value = 6243
$ws.Range("C488").Value = 5.16113018989563
$ws.Range("B489").Value = 4878
$ws.Range("C489").Value = 3.132783889770508
$ws.Range("B490").Value = 6616
$ws.Range("C490").Value = 5.567851066589355
$ws.Range("B491").Value = 16098
$ws.Range("C491").Value = 25.46987152099609
$ws.Range("D491").Value = 1
$ws.Range("B492").Value = 4822
$ws.Range("C492").Value = 3.019051313400269
$ws.Range("B493").Value = 1233
$ws.Range("C493").Value = 0.5142612457275391
$ws.Range("D493").Value = 1
$ws.Range("B494").Value = 7739
$ws.Range("C494").Value = 6.655339241027832
$ws.Range("D494").Value = 1
$ws.Range("B495").Value = 4743
$ws.Range("C495").Value = 3.030831098556519
$ws.Range("D495").Value = 1
$ws.Range("B496").Value = 19653
$ws.Range("C496").Value = 37.17036819458008
$ws.Range("D496").Value = 1
$ws.Range("B497").Value = 16449
$ws.Range("C497").Value = 24.47938013076782
$ws.Range("B498").Value = 4157
$ws.Range("C498").Value = 2.33838415145874
$ws.Range("B499").Value = 3742
$ws.Range("C499").Value = 2.350618362426758
$ws.Range("B500").Value = 2564
$ws.Range("C500").Value = 1.209393978118896
$ws.Range("B501").Value = 8374
$ws.Range("C501").Value = 7.737031698226929

$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("C2").Value = 0.198
